$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update 想去人数 (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10543
$ws1.Range("F3").Value = 234
$ws1.Range("F5").Value = 662

# Sheet "全部类型" (fourth sheet) - same updates mirrored here
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10543
$ws4.Range("F3").Value = 234
$ws4.Range("F5").Value = 662
